# Update "想去人数" (number interested) values in the F column
# on the "展览" and "全部类型" worksheets, to match the new data
# generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Rows 3,4,5,7,9,10,12,15 receive the same new values on both sheets.
# Row 13 differs slightly between the two sheets pre-edit (505 vs 506)
# but both converge to the same new value (514).
$commonUpdates = @{
    3  = 45
    4  = 1012
    5  = 14
    7  = 2564
    9  = 1609
    10 = 75
    12 = 65
    13 = 514
    15 = 50
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $commonUpdates.Keys) {
        $ws.Range("F$row").Value = $commonUpdates[$row]
    }
}
